$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (Reviews / Sentiment) right after the existing
# rows 1-13, mirroring the style used by the other "real" rows in column B.
$ws.Range("A14").Value = "Việt Nam có 1 triệu bệnh nhân covid"
$ws.Range("B14").Value = "real"

# Reflect the scrolled viewport / new selection left by the author when
# they saved the file (scrolled down to row 10, cursor resting on H19).
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("H19").Select()
